# Cotações atualizadas - 2025-11-02
# Append a new row (row 59) of quotes to the worksheet, mirroring the
# formatting/style used by the existing rows (numeric date in column A,
# text values in columns B-E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 59

# Column A: serial date number (same style as the rows above it, s="2")
$ws.Cells.Item($newRow, 1).Value = 45963

# Columns B-E: text values (decimal comma formatted strings), stored as text
$ws.Cells.Item($newRow, 2).Value = "22,0341"
$ws.Cells.Item($newRow, 3).Value = "16,1343"
$ws.Cells.Item($newRow, 4).Value = "15,5326"
$ws.Cells.Item($newRow, 5).Value = "15,5326"

# Match the style of the previous row's date cell (A58) for the new date cell (A59)
$ws.Range("A59").NumberFormat = $ws.Range("A58").NumberFormat
